{"js": "// Add a new \"Definition of Done\" bullet item right after the existing\n// \"Done means the feature has been developed, tested and meets all\n// required acceptance tests.\" bullet, using the same list (numId 9),\n// the same \"PlainText\" style, and the same Times New Roman / 22-half-point\n// run formatting as its neighbours.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that ends the \"Definition of Done\" bullet list.\nconst anchorText =\n  \"Done means the feature has been developed, tested and meets all required acceptance tests.\";\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the anchor paragraph for the new bullet.\");\n}\n\n// Insert a new paragraph right after it; insertParagraph copies the\n// paragraph/run formatting of the anchor (PlainText style, numPr numId=9,\n// Times New Roman 11pt runs), matching the surrounding bullets.\nconst newText =\n  \"Finish all process and meets the customer requirements then only it called done.\";\n\nconst newParagraph = anchor.insertParagraph(newText, Word.InsertLocation.after);\nnewParagraph.load(\"text\");\n\nawait context.sync();\n", "ps1": "# Add a new \"Definition of Done\" bullet item right after the existing\n# \"Done means the feature has been developed, tested and meets all\n# required acceptance tests.\" bullet. The new paragraph reuses the\n# same list (numId 9) / \"PlainText\" style / Times New Roman 11pt runs\n# as the rest of that bullet list, because it is created by splitting\n# the anchor paragraph, which already carries that formatting.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Done means the feature has been developed, tested and meets all required acceptance tests.\"\n\n# Locate the anchor paragraph with Find (robust against exact position\n# drift), then resolve the actual Paragraph object that contains it.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Could not find the anchor paragraph for the new bullet.\"\n}\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $pStart = $p.Range.Start\n    $pEnd = $p.Range.End\n    if ($pStart -le $rng.Start -and $pEnd -ge $rng.End) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not resolve the anchor Paragraph object.\"\n}\n\n# Split the anchor paragraph: insert a new empty paragraph right after\n# it (inherits pPr/rPr formatting), then fill in its text.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"Finish all process and meets the customer requirements then only it called done.\"\n"}
